$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.655.10'
$ws.Range('E2').Value = '  -7.37%  '
$ws.Range('D3').Value = '2.536.29'
$ws.Range('E3').Value = '  -3.53%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '298.49'
$ws.Range('E5').Value = '  -3.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '91.80'
$ws.Range('E6').Value = '  -7.00%  '
$ws.Range('E7').Value = '  -3.81%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.549'
$ws.Range('E9').Value = '  -5.35%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.07'
$ws.Range('E10').Value = '  -6.90%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0800'
$ws.Range('E11').Value = '  -5.41%  '
$ws.Range('E12').Value = '  -4.62%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.112'
$ws.Range('E13').Value = '  +4.68%  '
$ws.Range('D14').Value = '2.928.43'
$ws.Range('E14').Value = '  -3.34%  '
$ws.Range('D15').Value = '2.520.05'
$ws.Range('E15').Value = '  -4.32%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.872'
$ws.Range('E16').Value = '  -5.36%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.10'
$ws.Range('E17').Value = '  -5.27%  '
$ws.Range('D18').Value = '42.763.01'
$ws.Range('E18').Value = '  -7.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.00'
$ws.Range('E19').Value = '  +1.40%  '
$ws.Range('D20').Value = '0.0₃0980'
$ws.Range('E20').Value = '  -3.94%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.53'
$ws.Range('E21').Value = '  -3.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '71.45'
$ws.Range('E22').Value = '  -4.42%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '255.64'
$ws.Range('E23').Value = '  -9.85%  '
$ws.Range('E24').Value = '  -4.37%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '29.18'
$ws.Range('E25').Value = '  -2.87%  '
$ws.Range('E26').Value = '  -6.54%  '
$ws.Range('E27').Value = '  +0.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.00'
$ws.Range('E28').Value = '  -5.04%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.21'
$ws.Range('E29').Value = '  +0.20%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.82'
$ws.Range('E30').Value = '  -4.67%  '
$ws.Range('E31').Value = '  -4.93%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '151.70'
$ws.Range('E32').Value = '  -3.01%  '
$ws.Range('E33').Value = '  -6.34%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.75'
$ws.Range('E34').Value = '  -2.60%  '
$ws.Range('E35').Value = '  -8.15%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0792'
$ws.Range('E36').Value = '  -5.88%  '
$ws.Range('E37').Value = '  -7.46%  '
$ws.Range('B38').Value = 'Stellar'
$ws.Range('C38').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.119'
$ws.Range('E38').Value = '  -3.97%  '
$ws.Range('B39').Value = 'EnergySwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '24.07'
$ws.Range('E39').Value = '  +7.67%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '16.94'
$ws.Range('E40').Value = '  +7.02%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0309'
$ws.Range('E41').Value = '  -5.73%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.86'
$ws.Range('E42').Value = '  -4.24%  '
$ws.Range('E43').Value = '  -5.32%  '
$ws.Range('D44').Value = '2.079.92'
$ws.Range('E44').Value = '  -1.71%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.997'
$ws.Range('E45').Value = '  -0.19%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.07'
$ws.Range('E46').Value = '  -1.00%  '
$ws.Range('B47').Value = 'BitcoinSV'
$ws.Range('C47').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '84.32'
$ws.Range('E47').Value = '  -10.59%  '
$ws.Range('E48').Value = '  +2.24%  '
$ws.Range('D49').Value = '2.782.21'
$ws.Range('E49').Value = '  -3.35%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '103.99'
$ws.Range('E50').Value = '  -5.93%  '
$ws.Range('E51').Value = '  -4.99%  '
